$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.424.37"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "2.643.04"
$ws.Range("E3").Value = "  +1.99%  "
$ws.Range("E4").Value = "  -0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "589.20"
$c.NumberFormat = "General"
$ws.Range("E5").Value = "  +0.66%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "143.66"
$c.NumberFormat = "General"
$ws.Range("E6").Value = "  -2.36%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -1.62%  "
$ws.Range("D9").Value = "2.642.31"
$ws.Range("E9").Value = "  +2.03%  "
$ws.Range("E10").Value = "  -1.58%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.58"
$c.NumberFormat = "General"
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("E13").Value = "  -0.56%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "27.32"
$c.NumberFormat = "General"
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").Value = "3.117.72"
$ws.Range("E15").Value = "  +2.03%  "
$ws.Range("D16").Value = "63.360.74"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("E17").Value = "  -1.42%  "
$ws.Range("D18").Value = "2.676.72"
$ws.Range("E18").Value = "  +3.23%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.31"
$c.NumberFormat = "General"
$ws.Range("E19").Value = "  +0.04%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "339.96"
$c.NumberFormat = "General"
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("E21").Value = "  -1.30%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.68"
$c.NumberFormat = "General"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("E23").Value = "  +0.06%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "67.48"
$c.NumberFormat = "General"
$ws.Range("E24").Value = "  +0.36%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.65"
$c.NumberFormat = "General"
$ws.Range("E25").Value = "  +3.38%  "
$ws.Range("E26").Value = "  +5.88%  "
$ws.Range("E27").Value = "  -1.29%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "544.45"
$c.NumberFormat = "General"
$ws.Range("E28").Value = "  +13.95%  "
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("E31").Value = "  -1.64%  "
$ws.Range("E32").Value = "  +13.28%  "
$ws.Range("E33").Value = "  +1.63%  "
$ws.Range("D34").Value = "0.0₃0804"
$ws.Range("E34").Value = "  -2.49%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "173.76"
$c.NumberFormat = "General"
$ws.Range("E35").Value = "  -1.76%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  +6.57%  "
$ws.Range("E38").Value = "  -1.12%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "19.01"
$c.NumberFormat = "General"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("E40").Value = "  +3.73%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "170.39"
$c.NumberFormat = "General"
$ws.Range("E41").Value = "  +7.48%  "
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.NumberFormat = "General"
$ws.Range("E42").Value = "  -0.03%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "40.21"
$c.NumberFormat = "General"
$ws.Range("E43").Value = "  +1.87%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "3.72"
$c.NumberFormat = "General"
$ws.Range("E44").Value = "  -1.19%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "22.16"
$c.NumberFormat = "General"
$ws.Range("E45").Value = "  +3.76%  "
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("E47").Value = "  +1.17%  "
$ws.Range("E48").Value = "  -1.40%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "18.65"
$c.NumberFormat = "General"
$ws.Range("E50").Value = "  +1.64%  "
$ws.Range("E51").Value = "  -0.78%  "
